# Updated cryptos list (refreshed price/volume snapshot + two ranking swaps).
#
# All D/E cells in this sheet are stored as literal text (t="inlineStr"),
# e.g. "27.123.13" or "  +0.55%  " -- they are NOT real numbers/percentages.
# Excel's COM Range.Value setter auto-coerces any text that *parses* as a
# number (e.g. "1.007", "0.4686") into a true numeric value, which would
# change the cell's stored type/representation. To keep those cells as text
# (matching the source data), values that look numeric are written with a
# leading apostrophe, which is how Excel's "text entry" quote-prefix works.
# Values that aren't number-like (e.g. "27.108.66", names, URLs, the
# whitespace-padded percent strings) are written as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Cell,
        [string]$Text
    )
    if ($Text -match '^[+-]?\d+(\.\d+)?$') {
        # Looks like a plain number to Excel's auto-type-detection (e.g.
        # "1.007", "0.4686") -- prefix with an apostrophe so it is stored
        # as text (quote-prefixed), same as the source inline string.
        $ws.Range($Cell).Value = "'" + $Text
    } else {
        # Not number-like (multi-dot prices like "27.108.66", names, URLs,
        # the space-padded percent strings) -- safe to assign directly.
        $ws.Range($Cell).Value = $Text
    }
}

# row 2 - Bitcoin
Set-TextValue "D2" "27.108.66"
Set-TextValue "E2" "  +0.51%  "

# row 3 - Ethereum
Set-TextValue "D3" "1.825.96"
Set-TextValue "E3" "  +0.37%  "

# row 4 - TetherUSD
Set-TextValue "E4" "  +0.47%  "

# row 5 - BNB
Set-TextValue "D5" "312.63"
Set-TextValue "E5" "  +0.60%  "

# row 6 - USDC
Set-TextValue "D6" "1.007"
Set-TextValue "E6" "  +0.39%  "

# row 7 - XRP
Set-TextValue "D7" "0.4686"
Set-TextValue "E7" "  +0.26%  "

# row 8 - Cardano
Set-TextValue "E8" "  -0.37%  "

# row 9 - Dogecoin
Set-TextValue "D9" "0.07385"
Set-TextValue "E9" "  +0.50%  "

# row 10 - Polygon
Set-TextValue "D10" "0.8805"
Set-TextValue "E10" "  +0.77%  "

# row 11 - Solana
Set-TextValue "D11" "20.25"
Set-TextValue "E11" "  -0.21%  "

# row 12 - WrappedEther
Set-TextValue "D12" "1.849.05"
Set-TextValue "E12" "  +1.17%  "

# row 13 - TRON
Set-TextValue "D13" "0.07344"
Set-TextValue "E13" "  +2.94%  "

# row 14 - was Litecoin, now Polkadot (ranking swap with row 15)
Set-TextValue "B14" "Polkadot"
Set-TextValue "C14" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D14" "5.370"
Set-TextValue "E14" "  -0.87%  "

# row 15 - was Polkadot, now Litecoin (ranking swap with row 14)
Set-TextValue "B15" "Litecoin"
Set-TextValue "C15" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D15" "93.08"
Set-TextValue "E15" "  +1.66%  "

# row 16 - Chainlink
Set-TextValue "D16" "6.527"
Set-TextValue "E16" "  +0.16%  "

# row 17 - BinanceUSD
Set-TextValue "D17" "1.007"
Set-TextValue "E17" "  +0.21%  "

# row 18 - ShibaInu
Set-TextValue "D18" "0.000008720"
Set-TextValue "E18" "  -0.14%  "

# row 19 - Dai
Set-TextValue "D19" "1.007"
Set-TextValue "E19" "  +0.45%  "

# row 20 - WrappedBTC
Set-TextValue "D20" "27.462.90"
Set-TextValue "E20" "  +1.71%  "

# row 21 - Avalanche
Set-TextValue "D21" "14.62"

# row 22 - Uniswap
Set-TextValue "D22" "5.234"
Set-TextValue "E22" "  -1.16%  "

# row 23 - Cosmos
Set-TextValue "E23" "  -0.04%  "

# row 24 - WrappedliquidstakedEther2.0
Set-TextValue "D24" "2.080.81"
Set-TextValue "E24" "  +1.53%  "

# row 25 - Toncoin
Set-TextValue "D25" "1.882"
Set-TextValue "E25" "  -0.14%  "

# row 26 - Monero
Set-TextValue "D26" "151.16"
Set-TextValue "E26" "  -0.01%  "

# row 27 - EthereumClassic
Set-TextValue "D27" "18.49"
Set-TextValue "E27" "  +0.55%  "

# row 28 - LidoDAOToken
Set-TextValue "D28" "2.139"
Set-TextValue "E28" "  -0.09%  "

# row 29 - InternetComputer(DFINITY)
Set-TextValue "D29" "5.159"
Set-TextValue "E29" "  -1.79%  "

# row 30 - BitcoinCash
Set-TextValue "D30" "116.09"
Set-TextValue "E30" "  -0.76%  "

# row 31 - Stellar
Set-TextValue "D31" "0.08926"
Set-TextValue "E31" "  +0.47%  "

# row 32 - ImmutableX
Set-TextValue "D32" "0.7432"
Set-TextValue "E32" "  -1.89%  "

# row 33 - ARBITRUM
Set-TextValue "D33" "1.162"
Set-TextValue "E33" "  +0.20%  "

# row 34 - Filecoin
Set-TextValue "D34" "4.510"
Set-TextValue "E34" "  +0.13%  "

# row 35 - HuobiToken
Set-TextValue "D35" "2.943"
Set-TextValue "E35" "  +0.22%  "

# row 36 - Frax
Set-TextValue "D36" "1.007"
Set-TextValue "E36" "  +0.46%  "

# row 37 - was TrustWalletToken, now RenderToken (ranking swap with row 38)
Set-TextValue "B37" "RenderToken"
Set-TextValue "C37" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D37" "2.510"
Set-TextValue "E37" "  +5.45%  "

# row 38 - was RenderToken, now TrustWalletToken (ranking swap with row 37)
Set-TextValue "B38" "TrustWalletToken"
Set-TextValue "C38" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D38" "1.091"
Set-TextValue "E38" "  -0.75%  "

# row 39 - Hedera
Set-TextValue "D39" "0.05279"

# row 40 - VeChain
Set-TextValue "D40" "0.01937"
Set-TextValue "E40" "  -0.44%  "

# row 41 - FraxShare
Set-TextValue "D41" "7.336"
Set-TextValue "E41" "  +2.08%  "

# row 42 - MXToken
Set-TextValue "D42" "2.930"
Set-TextValue "E42" "  -1.40%  "

# row 43 - TheSandbox
Set-TextValue "D43" "0.5242"
Set-TextValue "E43" "  -1.14%  "

# row 44 - Algorand
Set-TextValue "E44" "  -0.74%  "

# row 45 - Aptos
Set-TextValue "D45" "8.384"
Set-TextValue "E45" "  -0.86%  "

# row 46 - Decentraland
Set-TextValue "D46" "0.4882"
Set-TextValue "E46" "  -0.21%  "

# row 47 - EnergySwap
Set-TextValue "D47" "10.41"
Set-TextValue "E47" "  -0.36%  "

# row 48 - PaxDollar
Set-TextValue "D48" "1.007"
Set-TextValue "E48" "  +0.46%  "

# row 49 - Quant
Set-TextValue "D49" "104.53"
Set-TextValue "E49" "  +1.19%  "

# row 50 - NEARProtocol
Set-TextValue "D50" "1.647"
Set-TextValue "E50" "  -1.06%  "
